$d = $word.ActiveDocument

# Match mode constants used with Find.Execute positional args:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#           MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Replace: 0=wdReplaceNone, 1=wdReplaceOne, 2=wdReplaceAll

# ---------------------------------------------------------------------------
# 1) Title "MusicPlayer" (144pt) - runs "Music"+"P"+"lay"+"er" -> single run
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("MusicPlayer", $true, $false, $false, $false, $false, $true, 1, $false, "MusicPlayer", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Heading "关于MusicPlayer" (40pt, bold) - merge runs
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("关于MusicPlayer", $true, $false, $false, $false, $false, $true, 1, $false, "关于MusicPlayer", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "实现技术" paragraph - merge "MusicPlayer" runs with the following
#    descriptive sentence run into one run
# ---------------------------------------------------------------------------
$rng = $d.Content
$text3 = "MusicPlayer是一个音乐播放器，在QtCreator环境下开发，使用了c++、qml和V-Play的插件来实现， 使用MYSQL数据库存储数据。"
$rng.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "功能简介" paragraph - merge "MusicPlayer" runs with the first part of
#    the descriptive sentence, then insert a new sentence in the middle as
#    its own run, before the trailing "支持两个用户..." run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$text4a = "MusicPlayer拥有音乐播放的相关功能，用户可以添加喜欢的音乐到指定歌单，允许用户自定义歌单、删除自定义歌单，以及向歌单里面添加和删除歌曲，"
$rng.Find.Execute($text4a, $true, $false, $false, $false, $false, $true, 1, $false, $text4a, 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute($text4a, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $rng.Duplicate
$insertPoint.Collapse(0)
$newText = "并将歌单信息保存在数据库中，"
$insertPoint.InsertAfter($newText)
$newRunRange = $d.Range($rng.End, $rng.End + $newText.Length)
# Force the freshly inserted text to become its own run (same visible
# formatting as its neighbours) by toggling a property and back again.
$newRunRange.Italic = 1
$newRunRange.Italic = 0

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from the "完善功能" paragraph to the start
#    of the "小组分工" heading paragraph.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$rng = $d.Content
$rng.Find.Execute("小组分工", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmTarget = $rng.Duplicate
$bmTarget.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmTarget) | Out-Null

# ---------------------------------------------------------------------------
# 6) "实现代码的合并" -> "实现歌单在数据库中的存取和读取", then add a new
#    paragraph right after it containing "实现代码的整合" (as two runs).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("实现代码的合并", $true, $false, $false, $false, $false, $true, 1, $false, "实现歌单在数据库中的存取和读取", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("实现歌单在数据库中的存取和读取", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphAfter()
$newParaStart = $rng.End + 1
$newParaRange = $d.Range($newParaStart, $newParaStart)
$newParaRange.InsertAfter("实现代码的")
$part2Start = $newParaStart + 5
$part2Range = $d.Range($part2Start, $part2Start)
$part2Range.InsertAfter("整合")
$part2RangeFmt = $d.Range($part2Start, $part2Start + 2)
$part2RangeFmt.Italic = 1
$part2RangeFmt.Italic = 0

# ---------------------------------------------------------------------------
# 7) "设计实现" section sentence - merge "Music"+"Pl"+"ay"+"er" runs plus
#    trailing text into one run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$text7 = "在设计实现阶段，实现了之前计划实现的功能，MusicPlayer基本功能实现。"
$rng.Find.Execute($text7, $true, $false, $false, $false, $false, $true, 1, $false, $text7, 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) "完善功能" section sentence - merge the trailing "MusicPlayer" run and
#    the "。" run (previously separated by the bookmark, now removed) into a
#    single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$text8 = "在完善功能阶段，我们经过讨论之后，确定了几个可新添加的功能，并分配了任务，最终完成MusicPlayer。"
$rng.Find.Execute($text8, $true, $false, $false, $false, $false, $true, 1, $false, $text8, 2) | Out-Null

Write-Host "Edits applied"
